$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. The "_GoBack" bookmark currently sits (collapsed) right after the
#    "draft project plan..." paragraph under Discussion. It needs to
#    move to sit between the "C" and "omber" of "Will comber" (which is
#    also being capitalised to "Will Comber" / "Will Comber"->split).
#    Remove it from its old spot first.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Capitalise "Will comber" -> "Will Comber" (fixing the attendee's
#    surname capitalisation).
# ------------------------------------------------------------------
$null = $d.Content.Find.Execute("Will comber", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "Will Comber", 2)

# ------------------------------------------------------------------
# 3. Locate that paragraph again (now reading "Will Comber") so we can
#    split it into three runs: "Will " / "C" / "omber", with the
#    "_GoBack" bookmark re-inserted between "C" and "omber".
# ------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Will Comber*") {
        $target = $p.Range
        break
    }
}

if ($target -ne $null) {
    $s = $target.Start

    # A collapsed range right after "Will " (before "C") forces a run
    # boundary there too, even though no bookmark should remain at
    # that spot -- add one temporarily, then remove it once the real
    # bookmark (which creates the boundary after "C") is in place.
    $splitPoint = $d.Range($s + 5, $s + 5)
    $d.Bookmarks.Add("TempRunSplit", $splitPoint)

    $goBackPoint = $d.Range($s + 6, $s + 6)
    $d.Bookmarks.Add("_GoBack", $goBackPoint)

    if ($d.Bookmarks.Exists("TempRunSplit")) {
        $d.Bookmarks("TempRunSplit").Delete()
    }
}
